# Apply updates to column F (dSF) for several rows, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value  = 0
$ws.Range("F5").Value  = 3
$ws.Range("F8").Value  = -1
$ws.Range("F15").Value = -4
$ws.Range("F21").Value = -1
$ws.Range("F22").Value = -5
$ws.Range("F24").Value = 3
$ws.Range("F25").Value = -8
